# Edit script implementing the "Small changes to compiled graphs" commit.
#
# 1) Slide 2 ("Advanced normalisation" slide): group the loose Picture 5 +
#    TextBox 14 shapes into a new Group so they move/scale together.
# 2) Slide 6 ("References" slide):
#    - Reposition/resize + embiggen+embolden the "References" title.
#    - Add two new reference textboxes (GANs bibliography, and a note about
#      the course the slide deck is based on).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: group Picture 5 (shape 5) + TextBox 14 (shape 6) -> Group 1
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$picShape = $null
$textShape = $null
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.Name -eq "Picture 5") { $picShape = $sh }
    if ($sh.Name -eq "TextBox 14") { $textShape = $sh }
}

$range2 = $s2.Shapes.Range(@($picShape.Id, $textShape.Id))
$grp = $range2.Group()

# ---------------------------------------------------------------------------
# Slide 6: References slide
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

$titleShape = $s6.Shapes.Item(1)   # "References"
$featureShape = $s6.Shapes.Item(2) # "Feature maps" bullet list

# Reposition + resize the title textbox, and bump up font size / bold.
$titleShape.Left = 395.0269291338583
$titleShape.Top = 22.36259842519685
$titleShape.Width = 261.72976377952756
$titleShape.Height = 36.3515748031496

$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Font.Size = 24
$titleRange.Font.Bold = 1

# New textbox: "GANs:" bibliography list.
$gansBox = $s6.Shapes.AddTextbox(1, 441.2432283464567, 82.39511811023623, 431.02700787401574, 378.05622047244094)
$gansBox.TextFrame.WordWrap = -1
$gansBox.TextFrame.AutoSize = 1

$gansText = "GANs:`rOverview of GANs: https://arxiv.org/pdf/1710.07035.pdf`rFace aging with cGAN: https://arxiv.org/pdf/1702.01983.pdf`rImage de-raining with cGAN: https://arxiv.org/pdf/1701.05957.pdf`rOriginal GAN paper: http://papers.nips.cc/paper/5423-generative-adversarial-nets.pdf`rDC-GAN paper: https://arxiv.org/pdf/1511.06434.pdf`rALI (Adversarially learned inference) model: https://arxiv.org/pdf/1606.00704.pdf`rBiGAN (same as ALI): https://arxiv.org/pdf/1605.09782.pdf`r`r"
$gansBox.TextFrame.TextRange.Text = $gansText

# Bullet formatting for paragraphs 2-10 (everything except the "GANs:" header).
for ($i = 2; $i -le $gansBox.TextFrame.TextRange.Paragraphs().Count; $i++) {
    $para = $gansBox.TextFrame.TextRange.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Visible = -1
    $para.ParagraphFormat.Bullet.Character = 8226
    $para.ParagraphFormat.Bullet.Font.Name = "Arial"
    $para.IndentLevel = 1
}

# Hyperlinks for each URL run within the GANs textbox.
$gansFull = $gansBox.TextFrame.TextRange.Text

function Add-Hyperlink($shapeTextRange, $fullText, $url) {
    $start = $fullText.IndexOf($url)
    if ($start -ge 0) {
        $sub = $shapeTextRange.Characters($start + 1, $url.Length)
        $sub.ActionSettings.Item(1).Hyperlink.Address = $url
    }
}

Add-Hyperlink $gansBox.TextFrame.TextRange $gansFull "https://arxiv.org/pdf/1710.07035.pdf"
Add-Hyperlink $gansBox.TextFrame.TextRange $gansFull "https://arxiv.org/pdf/1702.01983.pdf"
Add-Hyperlink $gansBox.TextFrame.TextRange $gansFull "https://arxiv.org/pdf/1701.05957.pdf"
Add-Hyperlink $gansBox.TextFrame.TextRange $gansFull "http://papers.nips.cc/paper/5423-generative-adversarial-nets.pdf"
Add-Hyperlink $gansBox.TextFrame.TextRange $gansFull "https://arxiv.org/pdf/1511.06434.pdf"
Add-Hyperlink $gansBox.TextFrame.TextRange $gansFull "https://arxiv.org/pdf/1606.00704.pdf"
Add-Hyperlink $gansBox.TextFrame.TextRange $gansFull "https://arxiv.org/pdf/1605.09782.pdf"

# New textbox: course notes reference.
$notesBox = $s6.Shapes.AddTextbox(1, 37.48795275590551, 428.291968503937, 431.02700787401574, 72.7031496062992)
$notesBox.TextFrame.WordWrap = -1
$notesBox.TextFrame.AutoSize = 1
$notesBox.TextFrame.TextRange.Text = "Notes on online deep learning course"
